$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 185, shifting existing rows 185-265 down to 186-266
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with the new record
$ws.Cells.Item(185, 1).Value  = 5
$ws.Cells.Item(185, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(185, 3).Value  = "Maule"
$ws.Cells.Item(185, 4).Value  = 44917
$ws.Cells.Item(185, 5).Value  = 7
$ws.Cells.Item(185, 6).Value  = 100112021
$ws.Cells.Item(185, 7).Value  = "Ají"
$ws.Cells.Item(185, 8).Value  = "Americana (o)"
$ws.Cells.Item(185, 9).Value  = "Primera"
$ws.Cells.Item(185, 10).Value = 100
$ws.Cells.Item(185, 11).Value = 13000
$ws.Cells.Item(185, 12).Value = 13000
$ws.Cells.Item(185, 13).Value = 13000
$ws.Cells.Item(185, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(185, 15).Value = "Región del Maule"
$ws.Cells.Item(185, 16).Value = 520
$ws.Cells.Item(185, 17).Value = 25
$ws.Cells.Item(185, 18).Value = "Hortaliza"
